$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that ends the work-log entry we need to extend:
# "Fix chức năng xác thực email." (a ListParagraph / bullet item).
# ------------------------------------------------------------------
$targetPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Fix chức năng xác thực email.*") {
        $targetPara = $p
    }
}

if ($targetPara -eq $null) {
    Write-Host "Could not find target paragraph"
}

# ------------------------------------------------------------------
# Reserve two fresh (empty) paragraphs right after the target one.
# They inherit the target's List Paragraph / numbering formatting,
# which is fine because we immediately overwrite each paragraph's
# full contents (pPr + run) with exact OOXML below.
# ------------------------------------------------------------------
$endRng = $targetPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$endRng.InsertParagraphAfter()

$targetIndex = $targetPara.Range.Information(3)
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- New paragraph 1: "12:03 AM 8/5 – Hiếu:" (plain Normal-style line) ---
$newPara1 = $targetPara.Next()
$xml1 = "<w:p $ns>" + `
    "<w:pPr>" + `
        "<w:tabs><w:tab w:val=`"left`" w:pos=`"864`"/></w:tabs>" + `
        "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:bCs/><w:sz w:val=`"26`"/><w:szCs w:val=`"26`"/></w:rPr>" + `
    "</w:pPr>" + `
    "<w:r>" + `
        "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:bCs/><w:sz w:val=`"26`"/><w:szCs w:val=`"26`"/></w:rPr>" + `
        "<w:lastRenderedPageBreak/>" + `
        "<w:t>12:03 AM 8/5 – Hiếu:</w:t>" + `
    "</w:r>" + `
    "</w:p>"
$fullRng1 = $d.Range($newPara1.Range.Start, $newPara1.Range.End)
$fullRng1.InsertXML($xml1)

# --- New paragraph 2: "Thêm trang admin" (bulleted ListParagraph item) ---
$newPara2 = $newPara1.Next()
$xml2 = "<w:p $ns>" + `
    "<w:pPr>" + `
        "<w:pStyle w:val=`"ListParagraph`"/>" + `
        "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr>" + `
        "<w:tabs><w:tab w:val=`"left`" w:pos=`"864`"/></w:tabs>" + `
        "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:bCs/><w:sz w:val=`"26`"/><w:szCs w:val=`"26`"/></w:rPr>" + `
    "</w:pPr>" + `
    "<w:r>" + `
        "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:bCs/><w:sz w:val=`"26`"/><w:szCs w:val=`"26`"/></w:rPr>" + `
        "<w:t>Thêm trang admin</w:t>" + `
    "</w:r>" + `
    "</w:p>"
$fullRng2 = $d.Range($newPara2.Range.Start, $newPara2.Range.End)
$fullRng2.InsertXML($xml2)

Write-Host "Inserted work-log entry for admin page setup."
